$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scheduled cryptos list refresh (GitHub Actions) - updated prices/volumes
# and a rotation of rows 12-15 (WrappedEther moved to rank 10).

$ws.Range('D2').Value = '26.418.06'
$ws.Range('D3').Value = '1.692.57'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').Value = "'219.01"
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = "'0.5539"
$ws.Range('E6').Value = '  +8.67%  '
$ws.Range('D7').Value = "'1.011"
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('D8').Value = "'0.2718"
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('D9').Value = "'0.06488"
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').Value = "'22.15"
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('D11').Value = "'0.07646"
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.692.11'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.563"
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.5825"
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = "'0.000008470"
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = "'65.28"
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('D17').Value = '26.494.72'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = "'4.962"
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = "'1.010"
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').Value = "'10.99"
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').Value = "'190.47"
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').Value = "'6.256"
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D24').Value = "'150.01"
$ws.Range('E24').Value = '  +3.35%  '
$ws.Range('E25').Value = '  +7.66%  '
$ws.Range('D26').Value = "'7.909"
$ws.Range('E26').Value = '  +3.90%  '
$ws.Range('D27').Value = "'15.78"
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').Value = "'1.427"
$ws.Range('E28').Value = '  +7.39%  '
$ws.Range('D29').Value = "'0.06336"
$ws.Range('E29').Value = '  -4.16%  '
$ws.Range('D30').Value = "'1.331"
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = "'3.599"
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').Value = "'3.598"
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').Value = "'1.679"
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('D34').Value = "'1.044"
$ws.Range('E34').Value = '  +2.71%  '
$ws.Range('D35').Value = "'0.6220"
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').Value = "'2.405"
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('D38').Value = "'6.233"
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '1.124.52'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('D40').Value = "'0.01644"
$ws.Range('E40').Value = '  +3.29%  '
$ws.Range('D41').Value = "'0.8816"
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('D43').Value = "'100.84"
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D44').Value = '1.844.18'
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('D46').Value = "'57.56"
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').Value = "'8.238"
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = "'0.05288"
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').Value = "'0.4304"
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').Value = "'6.082"
$ws.Range('E51').Value = '  +1.61%  '
